$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C gets wider to fit the "Blank hints" list; row 3 gets tall enough for the story text.
$ws.Columns.Item(3).ColumnWidth = 254.8
$ws.Rows.Item(3).RowHeight = 180

# Fill the hints column (C) first, then the story column (B), so the shared-string
# table picks up the same index order used by the target workbook (hints=5, story=6).
$ws.Range("C3").Value = @"
Adjective;Adjective;Celebrity;Adjective;Time frame - plural;Noun - Plural;Adjective - Ends in ER;Verb - Base Form;Letter of Alphabet;Part of Body;Name that starts with S;Color;Color;Celebrity;Verb - Present ends in S;Noun;Color;Number;Number;Month;Number;Year;Month;Number;Year;Number;Day of the week
"@

$ws.Range("B3").Value = @"
1. You know it's cold outside when you go outside, and it is {1}. 
2. {2} had been {3} {4} before she was dead.
3. Some {5} are {6} than others. 
4. The future has yet to {7}. 
5. The " {8} on Superman's {9} stands for {10}. 
6. The sky is {11}, and the grass is {12}. 
7. {13}'s last name is Obama. 
8. Rain {14} from the {15}. 
9. {16} is a color. 
10. Take your age, subtract {17}, then add {18}. That is your age. 
11. The 1950s lasted from {19} {20}, {21} - {22} {23}, {24}. 
12. Friday the 13th occurs when the {25}th day of a month falls on a {26}.
"@

$ws.Range("B3").WrapText = $true

[void]$ws.Range("C3").Select()
